$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestPlan")

# Row 2: switch unit-of-measure columns from KG to G (Gram), and rounding
# precision from 0.01 to 0.001
$ws.Range("F2").Value = "G, Gram, G"
$ws.Range("G2").Value = "G, Gram, G"
$ws.Range("H2").Value = "G, Gram, G"
$ws.Range("I2").Value = "0.001"
$ws.Range("J2").Value = "G, Gram, G"

# Row 3: rounding precision reverts back to 0.01
$ws.Range("I3").Value = "0.01"

# Update the active selection to I3 as per the saved view state
$ws.Activate()
$ws.Range("I3").Select()
